# Added output for tables I, III, IV.
#
# Table I's result columns are regenerated: the old "ols" / "within"
# estimator blocks are dropped, and the "infeasible" / "interactive"
# blocks are expanded from 2 coefficients + 2 std.devs each to 5 + 5
# each. Concretely, the header row (row 1) grows from columns B:S
# (i, t, coef1_ols, coef2_ols, sd1_ols, sd2_ols, coef1_within,
# coef2_within, sd1_within, sd2_within, coef1_infeasible,
# coef2_infeasible, sd1_infeasible, sd2_infeasible, coef1_interactive,
# coef2_interactive, sd1_interactive, sd2_interactive) to columns B:W
# (i, t, coef1..5_infeasible, sd1..5_infeasible, coef1..5_interactive,
# sd1..5_interactive).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = [ordered]@{
    "D1" = "coef1_infeasible"
    "E1" = "coef2_infeasible"
    "F1" = "coef3_infeasible"
    "G1" = "coef4_infeasible"
    "H1" = "coef5_infeasible"
    "I1" = "sd1_infeasible"
    "J1" = "sd2_infeasible"
    "K1" = "sd3_infeasible"
    "L1" = "sd4_infeasible"
    "M1" = "sd5_infeasible"
    "N1" = "coef1_interactive"
    "O1" = "coef2_interactive"
    "P1" = "coef3_interactive"
    "Q1" = "coef4_interactive"
    "R1" = "coef5_interactive"
    "S1" = "sd1_interactive"
    "T1" = "sd2_interactive"
    "U1" = "sd3_interactive"
    "V1" = "sd4_interactive"
    "W1" = "sd5_interactive"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}
